$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "82÷5=16, 2" -> "25÷5=5, 0"
$cell = $t.Cell(1, 1)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "25÷5=5, 0"

# Row 1, Col 2: "98÷8=12, 2" -> "87÷3=29, 0"
$cell = $t.Cell(1, 2)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "87÷3=29, 0"

# Row 1, Col 3: "23÷2=11, 1" -> "79÷2=39, 1"
$cell = $t.Cell(1, 3)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "79÷2=39, 1"

# Row 1, Col 4: "48÷8=6, 0" -> "93÷4=23, 1"
$cell = $t.Cell(1, 4)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "93÷4=23, 1"

# Row 1, Col 5: "76÷7=10, 6" -> "34÷6=5, 4"
$cell = $t.Cell(1, 5)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "34÷6=5, 4"

# Row 5, Col 1: "38÷4=9, 2" -> "47÷9=5, 2"
$cell = $t.Cell(5, 1)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "47÷9=5, 2"

# Row 5, Col 2: "51÷2=25, 1" -> "14÷2=7, 0"
$cell = $t.Cell(5, 2)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "14÷2=7, 0"

# Row 5, Col 3: "51÷8=6, 3" -> "43÷8=5, 3"
$cell = $t.Cell(5, 3)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "43÷8=5, 3"

# Row 5, Col 4: "83÷4=20, 3" -> "57÷6=9, 3"
$cell = $t.Cell(5, 4)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "57÷6=9, 3"

# Row 5, Col 5: "85÷9=9, 4" -> "65÷5=13, 0"
$cell = $t.Cell(5, 5)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "65÷5=13, 0"

# Row 9, Col 1: "50÷4=12, 2" -> "60÷2=30, 0"
$cell = $t.Cell(9, 1)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "60÷2=30, 0"

# Row 9, Col 2: "79÷9=8, 7" -> "58÷5=11, 3"
$cell = $t.Cell(9, 2)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "58÷5=11, 3"

# Row 9, Col 3: "78÷8=9, 6" -> "58÷3=19, 1"
$cell = $t.Cell(9, 3)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "58÷3=19, 1"

# Row 9, Col 4: "78÷4=19, 2" -> "83÷7=11, 6"
$cell = $t.Cell(9, 4)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "83÷7=11, 6"

# Row 9, Col 5: "78÷8=9, 6" -> "36÷3=12, 0"
$cell = $t.Cell(9, 5)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "36÷3=12, 0"

# Row 13, Col 1: "63÷5=12, 3" -> "68÷7=9, 5"
$cell = $t.Cell(13, 1)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "68÷7=9, 5"

# Row 13, Col 2: "39÷5=7, 4" -> "59÷2=29, 1"
$cell = $t.Cell(13, 2)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "59÷2=29, 1"

# Row 13, Col 3: "98÷2=49, 0" -> "55÷9=6, 1"
$cell = $t.Cell(13, 3)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "55÷9=6, 1"

# Row 13, Col 4: "98÷7=14, 0" -> "88÷4=22, 0"
$cell = $t.Cell(13, 4)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "88÷4=22, 0"

# Row 13, Col 5: "23÷5=4, 3" -> "56÷9=6, 2"
$cell = $t.Cell(13, 5)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "56÷9=6, 2"

# Row 17, Col 1: "69÷2=34, 1" -> "81÷6=13, 3"
$cell = $t.Cell(17, 1)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "81÷6=13, 3"

# Row 17, Col 2: "71÷5=14, 1" -> "99÷7=14, 1"
$cell = $t.Cell(17, 2)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "99÷7=14, 1"

# Row 17, Col 3: "92÷8=11, 4" -> "14÷9=1, 5"
$cell = $t.Cell(17, 3)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "14÷9=1, 5"

# Row 17, Col 4: "81÷7=11, 4" -> "10÷2=5, 0"
$cell = $t.Cell(17, 4)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "10÷2=5, 0"

# Row 17, Col 5: "62÷9=6, 8" -> "87÷2=43, 1"
$cell = $t.Cell(17, 5)
$r = $cell.Range
[void]$r.MoveEnd(1, -1)
$r.Text = "87÷2=43, 1"
